$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 93 with a new time-log entry (previously a blank template row).
# D93 is set before B93/C93 so the shared "=(C-B)*24-D/60" formula in E93
# picks up the interruption minutes on its first recalculation.
$ws.Range("A93").Value = 41931
$ws.Range("D93").Value = 20
$ws.Range("B93").Value = 0.7319444444444444
$ws.Range("C93").Value = 0.82361111111111107
$ws.Range("F93").Value = "Coding"

# Move the active selection to A94, matching the author's next-entry cursor
$ws.Range("A94").Select()
